# Updated AR TCs and Test Plan
#
# - "AddLine" sheet: remove the obsolete test-case row (old row 3), shifting
#   the remaining rows up and shrinking the used range from N6 to N5.
# - Make "AddLine" the active sheet/tab (it was "ARATO" before) and leave the
#   cursor on F7, matching the new selection captured in the sheet view.

$wb = $excel.ActiveWorkbook

$addLine = $wb.Worksheets.Item("AddLine")
$addLine.Rows.Item(3).EntireRow.Delete()

$addLine.Activate()
$addLine.Range("F7").Select()
